# Insert a new data row for "Terminal Hortofrutícola Agro Chillán - Pepino ensalada"
# right before the existing row 313, shifting all following rows down by one
# (old row 313 -> 314, ..., old row 372 -> 373), and fill the new row 313 with
# its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 313 (pushes rows 313:372 down to 314:373)
$ws.Rows.Item(313).Insert()

# Populate the newly inserted row 313
$ws.Cells.Item(313, 1).Value  = 7
$ws.Cells.Item(313, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(313, 3).Value  = "Ñuble"
$ws.Cells.Item(313, 4).Value  = 45173
$ws.Cells.Item(313, 5).Value  = 16
$ws.Cells.Item(313, 6).Value  = 100112043
$ws.Cells.Item(313, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(313, 8).Value  = "Sin especificar"
$ws.Cells.Item(313, 9).Value  = "Primera"
$ws.Cells.Item(313, 10).Value = 100
$ws.Cells.Item(313, 11).Value = 10000
$ws.Cells.Item(313, 12).Value = 10000
$ws.Cells.Item(313, 13).Value = 10000
$ws.Cells.Item(313, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(313, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(313, 16).Value = 167
$ws.Cells.Item(313, 17).Value = 60
$ws.Cells.Item(313, 18).Value = "Hortaliza"
